# Atualização automática via cronjob
# Adds the new 2025-05-09 sales-anomaly rows (12-15) and refreshes a handful
# of pre-existing values (A4, G4, A5, A6, G7, A8, G8, G9, A10, A11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Small in-place corrections to the existing rows (2-11)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = 8
$ws.Range("G4").Value = 339

$ws.Range("A5").Value = 11

$ws.Range("A6").Value = 13

$ws.Range("G7").Value = 865

$ws.Range("A8").Value = 5
$ws.Range("G8").Value = 184

$ws.Range("G9").Value = 12

$ws.Range("A10").Value = 9

$ws.Range("A11").Value = 12

# ---------------------------------------------------------------------
# 2) Append four new data rows (12-15), matching the style already used
#    for column A (bold / bordered / centered - same as A2:A11).
# ---------------------------------------------------------------------
$ws.Range("A11").Copy()
$ws.Range("A12:A15").PasteSpecial(-4122)

# -- Row 12 ------------------------------------------------------------
$ws.Range("A12").Value = 4
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "2025-05-09"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = "AMAZONIA REFEICOES E SERVICOS LTDA"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "000425"
$ws.Range("F12").Value = "COADOR DE CAFE EG (EXTRA GRANDE)"
$ws.Range("G12").Value = -3
$ws.Range("H12").Value = $false

# -- Row 13 ------------------------------------------------------------
$ws.Range("A13").Value = 6
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2025-05-09"
$ws.Range("C13").Value = 15
$ws.Range("D13").Value = "AMAZONIA REFEICOES E SERVICOS LTDA"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "001023"
$ws.Range("F13").Value = "FILME PVC 30X8X500MT"
$ws.Range("G13").Value = -1
$ws.Range("H13").Value = $true

# -- Row 14 ------------------------------------------------------------
$ws.Range("A14").Value = 7
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2025-05-09"
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = "AMAZONIA REFEICOES E SERVICOS LTDA"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "000924"
$ws.Range("F14").Value = "COPO POTE DESCARTAVEL TRANSP 100ML CX/20"
$ws.Range("G14").Value = -10
$ws.Range("H14").Value = $false

# -- Row 15 ------------------------------------------------------------
$ws.Range("A15").Value = 10
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "2025-05-09"
$ws.Range("C15").Value = 20
$ws.Range("D15").Value = "A R C DOS SANTOS"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "000770"
$ws.Range("F15").Value = "DETERGENTE CLORADO AUDAX GOLD 5L"
$ws.Range("G15").Value = 33
$ws.Range("H15").Value = $true

# ---------------------------------------------------------------------
# 3) Drop the Text number-format we used as a trick to stop Excel from
#    re-interpreting the date-like / zero-padded strings as numbers -
#    restores the cells to the workbook's default (un-styled) look.
# ---------------------------------------------------------------------
$ws.Range("B12:B15").ClearFormats()
$ws.Range("E12:E15").ClearFormats()

Write-Output "done"
